$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.782.92'
$ws.Range("E2").Value = '  -1.04%  '

$ws.Range("D3").Value = '2.601.97'
$ws.Range("E3").Value = '  -0.24%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.45'
$ws.Range("E5").Value = '  -1.22%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.69'
$ws.Range("E6").Value = '  -1.90%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.594'
$ws.Range("E7").Value = '  -0.97%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.575'
$ws.Range("E9").Value = '  -0.98%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.55'
$ws.Range("E10").Value = '  -0.56%  '

$ws.Range("E11").Value = '  -0.32%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0838'
$ws.Range("E12").Value = '  -0.49%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.03'
$ws.Range("E13").Value = '  -3.34%  '

$ws.Range("D14").Value = '2.997.50'
$ws.Range("E14").Value = '  -0.33%  '

$ws.Range("E15").Value = '  +0.78%  '

$ws.Range("D16").Value = '2.592.41'
$ws.Range("E16").Value = '  -0.58%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.909'
$ws.Range("E17").Value = '  +0.22%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.75'
$ws.Range("E18").Value = '  -2.05%  '

$ws.Range("D19").Value = '45.763.11'
$ws.Range("E19").Value = '  -1.49%  '

$ws.Range("E20").Value = '  -1.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.68'
$ws.Range("E21").Value = '  -1.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.65'
$ws.Range("E22").Value = '  -4.82%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '284.50'
$ws.Range("E23").Value = '  +11.64%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.62'
$ws.Range("E24").Value = '  +3.62%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.01'
$ws.Range("E25").Value = '  -1.98%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.24'
$ws.Range("E26").Value = '  +1.27%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '29.17'
$ws.Range("E27").Value = '  +3.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.16%  '

$ws.Range("E29").Value = '  +0.24%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.55'
$ws.Range("E30").Value = '  -0.19%  '

$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.25'
$ws.Range("E31").Value = '  -0.65%  '

$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '38.46'
$ws.Range("E32").Value = '  -5.08%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.22'
$ws.Range("E33").Value = '  +0.15%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.61'
$ws.Range("E34").Value = '  -3.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '157.11'
$ws.Range("E35").Value = '  +2.55%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.25'
$ws.Range("E36").Value = '  -1.52%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0830'
$ws.Range("E37").Value = '  -0.82%  '

$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.79'
$ws.Range("E38").Value = '  -3.23%  '

$ws.Range("E39").Value = '  +2.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.123'
$ws.Range("E40").Value = '  +0.14%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.81'
$ws.Range("E41").Value = '  -7.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0325'
$ws.Range("E42").Value = '  -0.57%  '

$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.52'
$ws.Range("E43").Value = '  -2.72%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.00'
$ws.Range("E44").Value = '  -5.61%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.28'
$ws.Range("E45").Value = '  +0.34%  '

$ws.Range("D46").Value = '2.099.16'
$ws.Range("E46").Value = '  +3.12%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.997'
$ws.Range("E47").Value = '  -0.26%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '93.82'
$ws.Range("E48").Value = '  +2.64%  '

$ws.Range("E49").Value = '  -0.96%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '108.62'
$ws.Range("E50").Value = '  -2.60%  '

$ws.Range("D51").Value = '2.848.75'
$ws.Range("E51").Value = '  -0.63%  '
